# Bug fix: remove empty footnotes - ones whose body carries no real
# content, just a lone Tibetan shad mark (punctuation) left over from a
# bad split. Iterate back-to-front so deleting doesn't renumber the
# items we still need to visit.
$d = $word.ActiveDocument

for ($i = $d.Footnotes.Count; $i -ge 1; $i--) {
    $fn = $d.Footnotes.Item($i)
    $text = $fn.Range.Text.Trim()
    if ($text.Length -le 1) {
        $fn.Delete()
    }
}
